# Applies the Monthly Report 5 text edits described by the commit diff.
# Each change is done as a narrow, unique Find/Replace so surrounding
# formatting (runs, paragraph properties) is preserved.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                                      $true, 1, $false, $new, 2)
    if (-not $found) {
        throw "Find/Replace failed for: $old"
    }
}

# Paragraph 1 ("With the end of my contract ...")
Replace-Text "end date of the project, fast approaching" `
             "end date of the first phase project, fast approaching"

Replace-Text "allow the passing through of the data. As we are using angular modals" `
             "allow the passing through of the data, next was getting it displayed on the screen. As we are using angular modals"

Replace-Text "could mainly be replicated with slight adjustments taking into account certain factors." `
             "could mainly be replicated with slight adjustments, taking into account certain factors."

# Paragraph 2 ("Seen as the project ...")
Replace-Text "consulting the division who it is being built form would be na" `
             "consulting the division, who it is being built for, would be na"

Replace-Text "where myself and the other would present the project from stretch to three people" `
             "where myself and the other intern along with our team leader and another team member would present the project from scratch to three people"

# Paragraph 3 ("With the nature of the schedule ...")
Replace-Text "With the nature of the schedule, this meeting was scheduled" `
             "Due to the nature of the respective schedules, this meeting was scheduled"

Replace-Text "At a meeting it was decided I would try build the application form as to give them a flavour" `
             "At a pre-meeting it was decided I would try build the application form in advance as to give them a flavour"

Replace-Text "Trying to create the application turned out to be difficult" `
             "Trying to create the application form turned out to be difficult"

# Paragraph 4 ("The meeting itself ...")
Replace-Text "When demoing the project everything was went over with a fine-combed by the members" `
             "When demoing the project, everything was examined over with a fine-comb by the members"
